$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "5001404-20.2017.8.21.0042"
$ws.Range("B2").Value = "9000985-29.2017.8.21.0042"
$ws.Range("C2").Value = "Migrado"

# Update row 3 values
$ws.Range("A3").Value = "5001392-06.2017.8.21.0042"
$ws.Range("B3").Value = "9001086-66.2017.8.21.0042"
$ws.Range("C3").Value = "Migrado"

# Remove rows 4-8 entirely (data no longer present)
$ws.Rows("4:8").Delete()
